# [DOCGEN] merge structure: plot + character lists
#
# Applies the template restructuring:
#   - Adds Author / Genre / Date written fields under the title
#   - Expands the single {{ character }} placeholder into a
#     {{ character_lists }} placeholder plus .name / .description /
#     .suggested_actors child placeholders (with a trailing blank line)
#   - Collapses the split "{{ " + "director"/"producer"/"writer" runs
#     back into single runs

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------
# 1) Author / Genre / Date written — inserted right after the title
#    paragraph, before the "Plot summary" heading.
# ---------------------------------------------------------------

$pTitle = $d.Paragraphs.Item(1)
$pTitle.Range.InsertParagraphAfter() | Out-Null
$pAuthor = $d.Paragraphs.Item(2)
$xmlAuthor = "<w:p $wNs>" +
  '<w:r><w:t xml:space="preserve">Author: </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>{{ author</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pAuthor.Range.InsertXML($xmlAuthor)

$pAuthor = $d.Paragraphs.Item(2)
$pAuthor.Range.InsertParagraphAfter() | Out-Null
$pGenre = $d.Paragraphs.Item(3)
$xmlGenre = "<w:p $wNs>" +
  '<w:r><w:t xml:space="preserve">Genre: </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>{{ genre</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pGenre.Range.InsertXML($xmlGenre)

$pGenre = $d.Paragraphs.Item(3)
$pGenre.Range.InsertParagraphAfter() | Out-Null
$pDateWritten = $d.Paragraphs.Item(4)
$xmlDateWritten = "<w:p $wNs>" +
  '<w:r><w:t>Date written: {{ date_written}}</w:t></w:r>' +
  '</w:p>'
$pDateWritten.Range.InsertXML($xmlDateWritten)

# ---------------------------------------------------------------
# 2) Character list: {{ character }} -> {{ character_lists }}, plus
#    .name / .description / .suggested_actors, plus a trailing blank
#    paragraph. The Author/Genre/Date-written inserts above shifted
#    everything that follows down by 3, so the old paragraph 5 is now
#    paragraph 8.
# ---------------------------------------------------------------

$pChar = $d.Paragraphs.Item(8)
$xmlChar = "<w:p $wNs>" +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t xml:space="preserve">{{ </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>character</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>_lists</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pChar.Range.InsertXML($xmlChar)

$pChar = $d.Paragraphs.Item(8)
$pChar.Range.InsertParagraphAfter() | Out-Null
$pName = $d.Paragraphs.Item(9)
$xmlName = "<w:p $wNs>" +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>{{ character_lists</w:t></w:r>' +
  '<w:r><w:t>.name</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pName.Range.InsertXML($xmlName)

$pName = $d.Paragraphs.Item(9)
$pName.Range.InsertParagraphAfter() | Out-Null
$pDescription = $d.Paragraphs.Item(10)
$xmlDescription = "<w:p $wNs>" +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t xml:space="preserve">{{ </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>character</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>_lists</w:t></w:r>' +
  '<w:r><w:t>.description</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pDescription.Range.InsertXML($xmlDescription)

$pDescription = $d.Paragraphs.Item(10)
$pDescription.Range.InsertParagraphAfter() | Out-Null
$pActors = $d.Paragraphs.Item(11)
$xmlActors = "<w:p $wNs>" +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t xml:space="preserve">{{ </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>character</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>_li</w:t></w:r>' +
  '<w:r><w:t>sts.suggested_actors</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pActors.Range.InsertXML($xmlActors)

$pActors = $d.Paragraphs.Item(11)
$pActors.Range.InsertParagraphAfter() | Out-Null
$pBlank = $d.Paragraphs.Item(12)
$xmlBlank = "<w:p $wNs></w:p>"
$pBlank.Range.InsertXML($xmlBlank)

# ---------------------------------------------------------------
# 3) Director / Producer / Writer: collapse the "{{ " + name (+ " }}")
#    run splits back into single runs. These paragraphs were shifted
#    down by 3 (Author/Genre/Date written) + 4 (character_lists
#    expansion: 3 new placeholder paragraphs + 1 blank paragraph).
# ---------------------------------------------------------------

$pDirector = $d.Paragraphs.Item(18)
$xmlDirector = "<w:p $wNs>" +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>{{ director</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pDirector.Range.InsertXML($xmlDirector)

$pProducer = $d.Paragraphs.Item(20)
$xmlProducer = "<w:p $wNs>" +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>{{ producer</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> }}</w:t></w:r>' +
  '</w:p>'
$pProducer.Range.InsertXML($xmlProducer)

$pWriter = $d.Paragraphs.Item(22)
$xmlWriter = "<w:p $wNs>" +
  '<w:r><w:t>{{ writer }}</w:t></w:r>' +
  '</w:p>'
$pWriter.Range.InsertXML($xmlWriter)
